$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 was an inline string "2"; change it to a real numeric value 2
$ws.Range("B4").Value = 2

# Add new row 5 with the additional annotation data
$ws.Range("A5").Value = "Ying Tang"
# Leading apostrophe forces this numeric-looking value to be stored as text
# (matches the source data, where politeness_score "1" is a text cell, not a number)
$ws.Range("B5").Value = "'1"
$ws.Range("C5").Value = "does not provide any insight, i tried this, i tried that ,a strong reject"
$ws.Range("D5").Value = "CRT"
$ws.Range("E5").Value = "THE"
$ws.Range("F5").Value = "1a2deef4-16ae-43c8-afd3-8fd2e076505e"
$ws.Range("G5").Value = "rJr4kfWCb_annotated.xlsx"
$ws.Range("H5").Value = "Overall, the paper does not provide any insight beyond: i tried this, i tried that and this works better than that; a strong reject."
